# The diff appends a period to the last sentence of the speaker notes on
# slide 25 ("...程序，算一下" -> "...程序，算一下。"), then adds a blank
# paragraph followed by a new parenthetical remark paragraph to the same
# notes text box.
#
# (Slide 25's notes explain the "measure without eBPF" baseline test; the
# new remark clarifies that this data set, together with the pictures on
# the next two slides - slide 26 "version 1" and slide 27 "version 2",
# each of which embeds exactly one picture - comes from a different test
# run than the numbers quoted in the final report.)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(25)
$np = $s.NotesPage

# The speaker-notes placeholder ("body idx=1") is the second shape on the
# notes page (1 = slide image placeholder, 2 = notes body, 3 = slide number).
$notesShape = $np.Shapes.Item(2)
$tr = $notesShape.TextFrame.TextRange

$nl = [char]10

$para1 = "那么接下来我们尝试测一下延迟：首先不挂eBPF程序，算一下。"
$para2 = ""
$para3 = "（本组数据与接下来的两个 ppt 中的图片，与结题报告结果分属两组不同时间的测试）"

$tr.Text = $para1 + $nl + $para2 + $nl + $para3
